# Apply "custom accuracy" rounding to row 5 values (round to 2 decimal places)
# and remove row 6 (the workbook is being trimmed as part of a larger,
# 1000-row dataset rebuild).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 measurement values down to 2 decimal places (custom accuracy).
$ws.Range("C5").Value = 7.96
$ws.Range("D5").Value = 0.92
$ws.Range("E5").Value = 24.01
$ws.Range("F5").Value = 19.26
$ws.Range("G5").Value = 8.6
$ws.Range("J5").Value = 5.85
$ws.Range("K5").Value = 8.539999999999999
$ws.Range("L5").Value = 9.619999999999999
$ws.Range("M5").Value = 10.21
$ws.Range("N5").Value = 2.78
$ws.Range("O5").Value = 8.65
$ws.Range("P5").Value = 12.23
$ws.Range("Q5").Value = 7.45
$ws.Range("R5").Value = 0.74
$ws.Range("S5").Value = 0.55
$ws.Range("T5").Value = 124.3
$ws.Range("U5").Value = 24.23
$ws.Range("V5").Value = 7.98
$ws.Range("W5").Value = 16.11
$ws.Range("X5").Value = 8.58
$ws.Range("Y5").Value = 1.46
$ws.Range("Z5").Value = 16.61
$ws.Range("AB5").Value = 6.32
$ws.Range("AD5").Value = 10.11
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 31.86
$ws.Range("AG5").Value = 4.41
$ws.Range("AH5").Value = 9.98

# Remove row 6 entirely so the used range shrinks to A1:AH5.
$ws.Rows.Item(6).Delete()
